# Auto-generated edit script: updates Leve market-price / profit columns
# (H..N) across multiple sheets to refresh cached Universalis market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 317.95834
$ws.Range("I28").Value = 340.95
$ws.Range("J28").Value = 203
$ws.Range("K28").Value = 340.95
$ws.Range("L28").Value = 203
$ws.Range("M28").Value = 144.05
$ws.Range("N28").Value = -1173

# Row 62
$ws.Range("H62").Value = 12809.182
$ws.Range("I62").Value = 22189.6
$ws.Range("J62").Value = 4992.1665
$ws.Range("K62").Value = 22189.6
$ws.Range("L62").Value = 4992.1665
$ws.Range("M62").Value = -21565.6
$ws.Range("N62").Value = -6240.1665

# Row 65
$ws.Range("H65").Value = 12809.182
$ws.Range("I65").Value = 22189.6
$ws.Range("J65").Value = 4992.1665
$ws.Range("K65").Value = 110948
$ws.Range("L65").Value = 24960.8325
$ws.Range("M65").Value = -107828
$ws.Range("N65").Value = -31200.8325

# Row 98
$ws.Range("H98").Value = 2347
$ws.Range("I98").Value = 1513.3572
$ws.Range("J98").Value = 5681.5713
$ws.Range("K98").Value = 1513.3572
$ws.Range("L98").Value = 5681.5713
$ws.Range("M98").Value = -15.35719999999992
$ws.Range("N98").Value = -8677.5713

# Row 101
$ws.Range("H101").Value = 1311
$ws.Range("I101").Value = 365.75
$ws.Range("J101").Value = 2571.3333
$ws.Range("K101").Value = 1097.25
$ws.Range("L101").Value = 7713.999899999999
$ws.Range("M101").Value = 524.75
$ws.Range("N101").Value = -10957.9999

# Row 116
$ws.Range("H116").Value = 64323.234
$ws.Range("I116").Value = 82955
$ws.Range("J116").Value = 3770
$ws.Range("K116").Value = 82955
$ws.Range("L116").Value = 3770
$ws.Range("M116").Value = -79513
$ws.Range("N116").Value = -10654

# Row 122
$ws.Range("H122").Value = 2347
$ws.Range("I122").Value = 1513.3572
$ws.Range("J122").Value = 5681.5713
$ws.Range("K122").Value = 4540.071599999999
$ws.Range("L122").Value = 17044.7139
$ws.Range("M122").Value = -2090.071599999999
$ws.Range("N122").Value = -21944.7139

# Row 132
$ws.Range("H132").Value = 4787.6274
$ws.Range("I132").Value = 1503.841
$ws.Range("J132").Value = 25428.572
$ws.Range("K132").Value = 4511.522999999999
$ws.Range("L132").Value = 76285.716
$ws.Range("M132").Value = -1981.522999999999
$ws.Range("N132").Value = -81345.716


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 785.1429000000001
$ws.Range("I2").Value = 719.2
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 719.2
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = -606.2
$ws.Range("N2").Value = -1176

# Row 32
$ws.Range("H32").Value = 5269.324
$ws.Range("I32").Value = 6342.6294
$ws.Range("J32").Value = 2371.4
$ws.Range("K32").Value = 6342.6294
$ws.Range("L32").Value = 2371.4
$ws.Range("M32").Value = -6055.6294
$ws.Range("N32").Value = -2945.4

# Row 45
$ws.Range("H45").Value = 1449.963
$ws.Range("I45").Value = 1221.75
$ws.Range("J45").Value = 2102
$ws.Range("K45").Value = 1221.75
$ws.Range("L45").Value = 2102
$ws.Range("M45").Value = -844.75
$ws.Range("N45").Value = -2856

# Row 61
$ws.Range("H61").Value = 718974.5
$ws.Range("I61").Value = 771762.7
$ws.Range("J61").Value = 673224.75
$ws.Range("K61").Value = 771762.7
$ws.Range("L61").Value = 673224.75
$ws.Range("M61").Value = -771550.7
$ws.Range("N61").Value = -673648.75

# Row 116
$ws.Range("H116").Value = 785.1429000000001
$ws.Range("I116").Value = 719.2
$ws.Range("J116").Value = 950
$ws.Range("K116").Value = 719.2
$ws.Range("L116").Value = 950
$ws.Range("M116").Value = 1574.8
$ws.Range("N116").Value = -5538

# Row 136
$ws.Range("H136").Value = 718974.5
$ws.Range("I136").Value = 771762.7
$ws.Range("J136").Value = 673224.75
$ws.Range("K136").Value = 2315288.1
$ws.Range("L136").Value = 2019674.25
$ws.Range("M136").Value = -2312738.1
$ws.Range("N136").Value = -2024774.25


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 785.1429000000001
$ws.Range("I3").Value = 719.2
$ws.Range("J3").Value = 950
$ws.Range("K3").Value = 719.2
$ws.Range("L3").Value = 950
$ws.Range("M3").Value = -605.2
$ws.Range("N3").Value = -1178


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1053.6666
$ws.Range("I16").Value = 1030.5
$ws.Range("K16").Value = 1030.5
$ws.Range("M16").Value = -743.5

# Row 31
$ws.Range("H31").Value = 2865.8667
$ws.Range("I31").Value = 992.8182
$ws.Range("J31").Value = 5155.148
$ws.Range("K31").Value = 992.8182
$ws.Range("L31").Value = 5155.148
$ws.Range("M31").Value = -697.8182
$ws.Range("N31").Value = -5745.148

# Row 34
$ws.Range("H34").Value = 2865.8667
$ws.Range("I34").Value = 992.8182
$ws.Range("J34").Value = 5155.148
$ws.Range("K34").Value = 992.8182
$ws.Range("L34").Value = 5155.148
$ws.Range("M34").Value = -790.8182
$ws.Range("N34").Value = -5559.148

# Row 99
$ws.Range("H99").Value = 49303.24
$ws.Range("I99").Value = 68053.2
$ws.Range("J99").Value = 2428.3333
$ws.Range("K99").Value = 68053.2
$ws.Range("L99").Value = 2428.3333
$ws.Range("M99").Value = -66555.2
$ws.Range("N99").Value = -5424.3333

# Row 113
$ws.Range("H113").Value = 1053.6666
$ws.Range("I113").Value = 1030.5
$ws.Range("K113").Value = 1030.5
$ws.Range("M113").Value = 1139.5

# Row 122
$ws.Range("H122").Value = 7066.6665
$ws.Range("I122").Value = 19800
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 59400
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -56950
$ws.Range("N122").Value = -7000

# Row 126
$ws.Range("H126").Value = 49303.24
$ws.Range("I126").Value = 68053.2
$ws.Range("J126").Value = 2428.3333
$ws.Range("K126").Value = 204159.6
$ws.Range("L126").Value = 7284.999899999999
$ws.Range("M126").Value = -201689.6
$ws.Range("N126").Value = -12224.9999

# Row 132
$ws.Range("H132").Value = 2253.7073
$ws.Range("I132").Value = 1555.1852
$ws.Range("J132").Value = 3600.8572
$ws.Range("K132").Value = 4665.5556
$ws.Range("L132").Value = 10802.5716
$ws.Range("M132").Value = -2135.5556
$ws.Range("N132").Value = -15862.5716


$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 797.9167
$ws.Range("I4").Value = 181.42857
$ws.Range("J4").Value = 1051.7646
$ws.Range("K4").Value = 544.28571
$ws.Range("L4").Value = 3155.2938
$ws.Range("M4").Value = -432.28571
$ws.Range("N4").Value = -3379.2938

# Row 5
$ws.Range("H5").Value = 732.5122
$ws.Range("I5").Value = 544.7143
$ws.Range("J5").Value = 929.7
$ws.Range("K5").Value = 1634.1429
$ws.Range("L5").Value = 2789.1
$ws.Range("M5").Value = -1522.1429
$ws.Range("N5").Value = -3013.1

# Row 12
$ws.Range("H12").Value = 278.5
$ws.Range("I12").Value = 197.25
$ws.Range("J12").Value = 319.125
$ws.Range("K12").Value = 591.75
$ws.Range("L12").Value = 957.375
$ws.Range("M12").Value = -418.75
$ws.Range("N12").Value = -1303.375

# Row 92
$ws.Range("H92").Value = 432.66666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 432.66666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1297.99998
$ws.Range("N92").Value = -3793.99998
$ws.Range("M92").ClearContents()

# Row 135
$ws.Range("H135").Value = 732.5122
$ws.Range("I135").Value = 544.7143
$ws.Range("J135").Value = 929.7
$ws.Range("K135").Value = 4902.428699999999
$ws.Range("L135").Value = 8367.300000000001
$ws.Range("M135").Value = -2367.428699999999
$ws.Range("N135").Value = -13437.3


$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 12502.5
$ws.Range("J5").Value = 12502.5
$ws.Range("L5").Value = 12502.5
$ws.Range("N5").Value = -12726.5

# Row 126
$ws.Range("H126").Value = 2545.3076
$ws.Range("I126").Value = 2677.8
$ws.Range("J126").Value = 2462.5
$ws.Range("K126").Value = 8033.400000000001
$ws.Range("L126").Value = 7387.5
$ws.Range("M126").Value = -5563.400000000001
$ws.Range("N126").Value = -12327.5


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2039.9706
$ws.Range("I40").Value = 2001.6333
$ws.Range("J40").Value = 2327.5
$ws.Range("K40").Value = 2001.6333
$ws.Range("L40").Value = 2327.5
$ws.Range("M40").Value = -1865.6333
$ws.Range("N40").Value = -1999.5

